$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.922.38"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "3.146.93"
$ws.Range("E3").Value = "  -7.49%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "568.48"
$ws.Range("E5").Value = "  -2.46%  "
$ws.Range("D6").Value = "169.16"
$ws.Range("E6").Value = "  -5.94%  "
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "3.142.98"
$ws.Range("E9").Value = "  -7.61%  "
$ws.Range("E10").Value = "  -5.36%  "
$ws.Range("D11").Value = "6.53"
$ws.Range("E11").Value = "  -6.20%  "
$ws.Range("D12").Value = "0.391"
$ws.Range("E12").Value = "  -5.02%  "
$ws.Range("D13").Value = "3.695.87"
$ws.Range("E13").Value = "  -7.44%  "
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "26.97"
$ws.Range("D16").Value = "64.802.02"
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("D17").Value = "0.0000161"
$ws.Range("E17").Value = "  -6.22%  "
$ws.Range("D18").Value = "3.156.08"
$ws.Range("E18").Value = "  -7.43%  "
$ws.Range("D19").Value = "5.71"
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("D20").Value = "12.76"
$ws.Range("E20").Value = "  -7.37%  "
$ws.Range("D21").Value = "357.03"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "7.24"
$ws.Range("E22").Value = "  -4.39%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "68.96"
$ws.Range("E24").Value = "  -5.58%  "
$ws.Range("D25").Value = "0.495"
$ws.Range("E25").Value = "  -6.74%  "
$ws.Range("D26").Value = "3.303.50"
$ws.Range("E26").Value = "  -7.45%  "
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  -7.40%  "
$ws.Range("D28").Value = "9.67"
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").Value = "1.90"
$ws.Range("E32").Value = "  -4.43%  "
$ws.Range("D33").Value = "21.90"
$ws.Range("E33").Value = "  -5.60%  "
$ws.Range("D34").Value = "5.30"
$ws.Range("E34").Value = "  -7.69%  "
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("D36").Value = "6.59"
$ws.Range("E36").Value = "  -6.06%  "
$ws.Range("D37").Value = "158.06"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("D38").Value = "1.43"
$ws.Range("E38").Value = "  -6.89%  "
$ws.Range("D39").Value = "0.832"
$ws.Range("E39").Value = "  -3.76%  "
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").Value = "25.99"
$ws.Range("E41").Value = "  -5.04%  "
$ws.Range("D42").Value = "2.666.08"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "2.44"
$ws.Range("D44").Value = "4.16"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("D45").Value = "6.03"
$ws.Range("E45").Value = "  -3.64%  "
$ws.Range("D46").Value = "39.43"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "0.0652"
$ws.Range("E47").Value = "  -4.64%  "
$ws.Range("D48").Value = "23.97"
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("D49").Value = "319.52"
$ws.Range("E49").Value = "  -3.94%  "
$ws.Range("D50").Value = "0.0271"
$ws.Range("E50").Value = "  -4.77%  "
$ws.Range("E51").Value = "  -1.34%  "
